{"js": "// The document contains a firewalld example command that ends with\n// \" \u2013permanent\" where the dash before \"permanent\" is a Unicode EN DASH\n// (U+2013) instead of two regular hyphens. Fix it to read \" --permanent\"\n// (matching the double-hyphen long-option style used by the rest of the\n// command, e.g. \"--zone=public\", \"--add-port\").\n\nconst searchText = \" \\u2013permanent\"; // \" \u2013permanent\"\nconst replacementText = \" --permanent\";\n\nconst results = context.document.body.search(searchText, { matchCase: true });\nresults.load(\"items,text\");\nawait context.sync();\n\nfor (let i = 0; i < results.items.length; i++) {\n  results.items[i].insertText(replacementText, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# The document contains a firewalld example command that ends with\n# \" -permanent\" where the dash before \"permanent\" is a Unicode EN DASH\n# (U+2013) instead of two regular hyphens. Fix it to read \" --permanent\"\n# (matching the double-hyphen long-option style used by the rest of the\n# command, e.g. \"--zone=public\", \"--add-port\").\n\n$d = $word.ActiveDocument\n\n$enDash = [char]0x2013\n$findText = \" \" + $enDash + \"permanent\"\n$replaceText = \" --permanent\"\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = $findText\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = $replaceText\n$find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)\n"}
